{"js": "// Replace each \"before\" math expression with its corresponding \"after\"\n// expression, one per table cell. Every \"before\" string is unique in the\n// document and matches exactly one run, so a search+replace per pair is\n// safe and order-independent.\nconst replacements = [\n  [\"82-58=24\", \"26-22=4\"],\n  [\"83-81=2\", \"87-51=36\"],\n  [\"85-80=5\", \"23+40=63\"],\n  [\"79-36=43\", \"57+28=85\"],\n  [\"8-6=2\", \"82+9=91\"],\n  [\"57-55=2\", \"25+31=56\"],\n  [\"63-58=5\", \"8+78=86\"],\n  [\"22+47=69\", \"72-26=46\"],\n  [\"62+20=82\", \"86-34=52\"],\n  [\"98-84=14\", \"52-38=14\"],\n  [\"86-40=46\", \"96-14=82\"],\n  [\"40-29=11\", \"28+68=96\"],\n  [\"2+72=74\", \"8+38=46\"],\n  [\"10-5=5\", \"95-90=5\"],\n  [\"94-93=1\", \"23-14=9\"],\n  [\"40+38=78\", \"55+12=67\"],\n  [\"4+66=70\", \"7+61=68\"],\n  [\"70-20=50\", \"55+15=70\"],\n  [\"42-18=24\", \"28+51=79\"],\n  [\"15+35=50\", \"38-20=18\"],\n  [\"80+0=80\", \"64+3=67\"],\n  [\"48-4=44\", \"14+62=76\"],\n  [\"90-7=83\", \"46-20=26\"],\n  [\"0+28=28\", \"42-24=18\"],\n  [\"41-21=20\", \"70-13=57\"],\n  [\"88+7=95\", \"56+12=68\"],\n  [\"17+63=80\", \"23-3=20\"],\n  [\"31-17=14\", \"39-25=14\"],\n  [\"37+57=94\", \"37+26=63\"],\n  [\"83-3=80\", \"84-59=25\"],\n  [\"82+7=89\", \"74-56=18\"],\n  [\"43+4=47\", \"85-32=53\"],\n  [\"31+2=33\", \"29-5=24\"],\n  [\"57-29=28\", \"32-11=21\"],\n  [\"99-41=58\", \"48-44=4\"],\n  [\"77-66=11\", \"84-81=3\"],\n  [\"29+38=67\", \"36+39=75\"],\n  [\"35+25=60\", \"55+3=58\"],\n  [\"67-25=42\", \"71-31=40\"],\n  [\"82-59=23\", \"40+50=90\"],\n  [\"12+42=54\", \"60+22=82\"],\n  [\"65-3=62\", \"7+59=66\"],\n  [\"89-30=59\", \"30+23=53\"],\n  [\"95-12=83\", \"10+33=43\"],\n  [\"89-19=70\", \"24+9=33\"],\n  [\"50-37=13\", \"30+17=47\"],\n  [\"46+18=64\", \"12+7=19\"],\n  [\"29+28=57\", \"17+39=56\"],\n  [\"40+28=68\", \"84-23=61\"],\n  [\"20-9=11\", \"35-34=1\"],\n  [\"58-36=22\", \"80-32=48\"],\n  [\"8+76=84\", \"28-12=16\"],\n  [\"2+0=2\", \"58+32=90\"],\n  [\"69+12=81\", \"39-10=29\"],\n  [\"0+61=61\", \"74+12=86\"],\n  [\"2+74=76\", \"88-1=87\"],\n  [\"47-8=39\", \"82-34=48\"],\n  [\"54-31=23\", \"0+10=10\"],\n  [\"36+63=99\", \"0+5=5\"],\n  [\"79+14=93\", \"68-34=34\"],\n  [\"41+40=81\", \"49-38=11\"],\n  [\"96-6=90\", \"7+18=25\"],\n  [\"27+47=74\", \"6+68=74\"],\n  [\"6+13=19\", \"82-9=73\"],\n  [\"81-58=23\", \"10-7=3\"],\n  [\"72-27=45\", \"38-13=25\"],\n  [\"41-11=30\", \"93-47=46\"],\n  [\"53-19=34\", \"53+23=76\"],\n  [\"68+29=97\", \"2+50=52\"],\n  [\"20+4=24\", \"23+34=57\"],\n  [\"87-35=52\", \"35+62=97\"],\n  [\"82-37=45\", \"93+6=99\"],\n  [\"82-82=0\", \"18+16=34\"],\n  [\"5+16=21\", \"40-11=29\"],\n  [\"26+1=27\", \"73-4=69\"],\n  [\"97-15=82\", \"68-63=5\"],\n  [\"4+37=41\", \"15+74=89\"],\n  [\"1+15=16\", \"2+65=67\"],\n  [\"30+35=65\", \"81-71=10\"],\n  [\"67-10=57\", \"11+45=56\"],\n  [\"18-8=10\", \"93-37=56\"],\n  [\"14-5=9\", \"30-14=16\"],\n  [\"89-84=5\", \"82-39=43\"],\n  [\"72-1=71\", \"99-67=32\"],\n  [\"83-43=40\", \"24-23=1\"],\n  [\"33+17=50\", \"16+30=46\"],\n  [\"56-34=22\", \"48+39=87\"],\n  [\"21+47=68\", \"18+25=43\"],\n  [\"60+15=75\", \"59-25=34\"],\n  [\"84-33=51\", \"44-25=19\"],\n  [\"63-26=37\", \"87-49=38\"],\n  [\"70-25=45\", \"79-31=48\"],\n  [\"96-47=49\", \"71-59=12\"],\n  [\"78-34=44\", \"36-11=25\"],\n  [\"11+22=33\", \"98-0=98\"],\n  [\"46+2=48\", \"40-27=13\"],\n  [\"16+12=28\", \"40-21=19\"],\n  [\"45+2=47\", \"50-10=40\"],\n  [\"5+19=24\", \"28-25=3\"],\n  [\"21+31=52\", \"40+58=98\"]\n];\n\nconst body = context.document.body;\n\nfor (const [before, after] of replacements) {\n  const results = body.search(before, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`No match found for \"${before}\"`);\n  }\n\n  results.items.forEach((item) => {\n    item.insertText(after, Word.InsertLocation.replace);\n  });\n}\n\nawait context.sync();\n", "ps1": "# Replace each \"before\" math expression with its corresponding \"after\"\n# expression, one per table cell. Every \"before\" string is unique in the\n# document (no duplicates, no substrings of each other), so a sequential\n# Find/Replace per pair is safe and order-independent.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"82-58=24\", \"26-22=4\"),\n    @(\"83-81=2\", \"87-51=36\"),\n    @(\"85-80=5\", \"23+40=63\"),\n    @(\"79-36=43\", \"57+28=85\"),\n    @(\"8-6=2\", \"82+9=91\"),\n    @(\"57-55=2\", \"25+31=56\"),\n    @(\"63-58=5\", \"8+78=86\"),\n    @(\"22+47=69\", \"72-26=46\"),\n    @(\"62+20=82\", \"86-34=52\"),\n    @(\"98-84=14\", \"52-38=14\"),\n    @(\"86-40=46\", \"96-14=82\"),\n    @(\"40-29=11\", \"28+68=96\"),\n    @(\"2+72=74\", \"8+38=46\"),\n    @(\"10-5=5\", \"95-90=5\"),\n    @(\"94-93=1\", \"23-14=9\"),\n    @(\"40+38=78\", \"55+12=67\"),\n    @(\"4+66=70\", \"7+61=68\"),\n    @(\"70-20=50\", \"55+15=70\"),\n    @(\"42-18=24\", \"28+51=79\"),\n    @(\"15+35=50\", \"38-20=18\"),\n    @(\"80+0=80\", \"64+3=67\"),\n    @(\"48-4=44\", \"14+62=76\"),\n    @(\"90-7=83\", \"46-20=26\"),\n    @(\"0+28=28\", \"42-24=18\"),\n    @(\"41-21=20\", \"70-13=57\"),\n    @(\"88+7=95\", \"56+12=68\"),\n    @(\"17+63=80\", \"23-3=20\"),\n    @(\"31-17=14\", \"39-25=14\"),\n    @(\"37+57=94\", \"37+26=63\"),\n    @(\"83-3=80\", \"84-59=25\"),\n    @(\"82+7=89\", \"74-56=18\"),\n    @(\"43+4=47\", \"85-32=53\"),\n    @(\"31+2=33\", \"29-5=24\"),\n    @(\"57-29=28\", \"32-11=21\"),\n    @(\"99-41=58\", \"48-44=4\"),\n    @(\"77-66=11\", \"84-81=3\"),\n    @(\"29+38=67\", \"36+39=75\"),\n    @(\"35+25=60\", \"55+3=58\"),\n    @(\"67-25=42\", \"71-31=40\"),\n    @(\"82-59=23\", \"40+50=90\"),\n    @(\"12+42=54\", \"60+22=82\"),\n    @(\"65-3=62\", \"7+59=66\"),\n    @(\"89-30=59\", \"30+23=53\"),\n    @(\"95-12=83\", \"10+33=43\"),\n    @(\"89-19=70\", \"24+9=33\"),\n    @(\"50-37=13\", \"30+17=47\"),\n    @(\"46+18=64\", \"12+7=19\"),\n    @(\"29+28=57\", \"17+39=56\"),\n    @(\"40+28=68\", \"84-23=61\"),\n    @(\"20-9=11\", \"35-34=1\"),\n    @(\"58-36=22\", \"80-32=48\"),\n    @(\"8+76=84\", \"28-12=16\"),\n    @(\"2+0=2\", \"58+32=90\"),\n    @(\"69+12=81\", \"39-10=29\"),\n    @(\"0+61=61\", \"74+12=86\"),\n    @(\"2+74=76\", \"88-1=87\"),\n    @(\"47-8=39\", \"82-34=48\"),\n    @(\"54-31=23\", \"0+10=10\"),\n    @(\"36+63=99\", \"0+5=5\"),\n    @(\"79+14=93\", \"68-34=34\"),\n    @(\"41+40=81\", \"49-38=11\"),\n    @(\"96-6=90\", \"7+18=25\"),\n    @(\"27+47=74\", \"6+68=74\"),\n    @(\"6+13=19\", \"82-9=73\"),\n    @(\"81-58=23\", \"10-7=3\"),\n    @(\"72-27=45\", \"38-13=25\"),\n    @(\"41-11=30\", \"93-47=46\"),\n    @(\"53-19=34\", \"53+23=76\"),\n    @(\"68+29=97\", \"2+50=52\"),\n    @(\"20+4=24\", \"23+34=57\"),\n    @(\"87-35=52\", \"35+62=97\"),\n    @(\"82-37=45\", \"93+6=99\"),\n    @(\"82-82=0\", \"18+16=34\"),\n    @(\"5+16=21\", \"40-11=29\"),\n    @(\"26+1=27\", \"73-4=69\"),\n    @(\"97-15=82\", \"68-63=5\"),\n    @(\"4+37=41\", \"15+74=89\"),\n    @(\"1+15=16\", \"2+65=67\"),\n    @(\"30+35=65\", \"81-71=10\"),\n    @(\"67-10=57\", \"11+45=56\"),\n    @(\"18-8=10\", \"93-37=56\"),\n    @(\"14-5=9\", \"30-14=16\"),\n    @(\"89-84=5\", \"82-39=43\"),\n    @(\"72-1=71\", \"99-67=32\"),\n    @(\"83-43=40\", \"24-23=1\"),\n    @(\"33+17=50\", \"16+30=46\"),\n    @(\"56-34=22\", \"48+39=87\"),\n    @(\"21+47=68\", \"18+25=43\"),\n    @(\"60+15=75\", \"59-25=34\"),\n    @(\"84-33=51\", \"44-25=19\"),\n    @(\"63-26=37\", \"87-49=38\"),\n    @(\"70-25=45\", \"79-31=48\"),\n    @(\"96-47=49\", \"71-59=12\"),\n    @(\"78-34=44\", \"36-11=25\"),\n    @(\"11+22=33\", \"98-0=98\"),\n    @(\"46+2=48\", \"40-27=13\"),\n    @(\"16+12=28\", \"40-21=19\"),\n    @(\"45+2=47\", \"50-10=40\"),\n    @(\"5+19=24\", \"28-25=3\"),\n    @(\"21+31=52\", \"40+58=98\")\n)\n\n# Find.Execute positional args:\n#   FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n#   MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace\n# Wrap=1 -> wdFindContinue, Replace=2 -> wdReplaceAll\nforeach ($pair in $pairs) {\n    $before = $pair[0]\n    $after = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n\n    $found = $find.Execute($before, $false, $false, $false, $false, $false, $true, 1, $false, $after, 2)\n    if (-not $found) {\n        Write-Output \"WARNING: not found -> $before\"\n    }\n}\n"}
